$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the now-unused trailing columns (X:AG) from the header row and the
#    second (label) row - the sample set shrank from 32 columns to 22.
# ---------------------------------------------------------------------------
$ws.Range("X1:AG2").Clear()

# ---------------------------------------------------------------------------
# 2. The simulation-scheme name list was overhauled. Relabel column B for
#    every existing data row (3-19) with the new scheme names...
# ---------------------------------------------------------------------------
$colB = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD"
)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(3 + $i, 2).Value = $colB[$i]
}

# ...and replace the scheme labels across row 2 (columns C:W) with the new
# [h,k,l]/pairing labels.
$row2 = @("[3, 3, 1]","[3, 1, 1]","[1, 1, 1]","[2, 2, 2]","[5, 1, 1]","[4, 2, 2]","[4, 2, 0]","[4, 0, 0]","[3, 3, 3]","[2, 0, 0]","[2, 2, 0]","1Pair-A","1Pair-B","2Pairs-A","2Pairs-B","3Pairs-A","3Pairs-B","3Pairs-C","4Pairs","5A4F","MaxUnique")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $row2[$i]
}

# ---------------------------------------------------------------------------
# 3. Append the new simulation rows (20-29), matching the existing row
#    layout/formatting (bold bordered index in column A, label in column B,
#    and a run of 1's across C:W).
# ---------------------------------------------------------------------------
$newRows = @(
    @(18, "Ring Perpendicular to TD"),
    @(19, "OffsetFTD"),
    @(20, "OffsetATD"),
    @(21, "OffsetF45"),
    @(22, "OffsetA45"),
    @(23, "OffsetFRD"),
    @(24, "OffsetARD"),
    @(25, "Gaussian Quadrature"),
    @(26, "Michael-CCHex"),
    @(27, "Michael-SNHex")
)

$destRow = 20
foreach ($entry in $newRows) {
    $ws.Range("A19:W19").Copy()
    $ws.Range("A" + $destRow + ":W" + $destRow).PasteSpecial(-4122)

    $ws.Cells.Item($destRow, 1).Value = $entry[0]
    $ws.Cells.Item($destRow, 2).Value = $entry[1]
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($destRow, $c).Value = 1
    }

    $destRow = $destRow + 1
}
